$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 6.9
$ws.Range("C3").Value = 6.2
$ws.Range("C11").Select()
